# This script updates the TPM-derived cell-cell communication metrics
# (NATMI lrc2p output) for the Col1a2-Cd44 ligand-receptor pair sheet,
# reflecting new TPM values used to recompute expression / specificity
# statistics in columns G-J (ligand), M-P (receptor), and Q-T (edge).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.042494666666667
$ws.Range("H2").Value = 3.127484
$ws.Range("I2").Value = 0.0007670466909205676
$ws.Range("J2").Value = 0.0007670466909205677
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 7.805806001015998
$ws.Range("R2").Value = 70.252254009144
$ws.Range("S2").Value = 0.0001061895741157796
$ws.Range("T2").Value = 0.0001061895741157796
$ws.Range("G3").Value = 1.042494666666667
$ws.Range("H3").Value = 3.127484
$ws.Range("I3").Value = 0.0007670466909205676
$ws.Range("J3").Value = 0.0007670466909205677
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 33.359233026384
$ws.Range("R3").Value = 300.233097237456
$ws.Range("S3").Value = 0.0004538163960825682
$ws.Range("T3").Value = 0.0004538163960825683
$ws.Range("G4").Value = 1.042494666666667
$ws.Range("H4").Value = 3.127484
$ws.Range("I4").Value = 0.0007670466909205676
$ws.Range("J4").Value = 0.0007670466909205677
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 15.21919372711778
$ws.Range("R4").Value = 136.97274354406
$ws.Range("S4").Value = 0.0002070407207222197
$ws.Range("T4").Value = 0.0002070407207222198
$ws.Range("I5").Value = 0.9658609009611662
$ws.Range("J5").Value = 0.9658609009611662
$ws.Range("M5").Value = 7.487621999999999
$ws.Range("N5").Value = 22.462866
$ws.Range("O5").Value = 0.1384395179233961
$ws.Range("P5").Value = 0.1384395179233961
$ws.Range("Q5").Value = 9829.027236687647
$ws.Range("R5").Value = 88461.24513018882
$ws.Range("S5").Value = 0.1337133175101208
$ws.Range("T5").Value = 0.1337133175101208
$ws.Range("I6").Value = 0.9658609009611662
$ws.Range("J6").Value = 0.9658609009611662
$ws.Range("O6").Value = 0.5916411627275552
$ws.Range("P6").Value = 0.5916411627275552
$ws.Range("S6").Value = 0.5714430664777483
$ws.Range("T6").Value = 0.5714430664777483
$ws.Range("I7").Value = 0.9658609009611662
$ws.Range("J7").Value = 0.9658609009611662
$ws.Range("M7").Value = 14.59882166666667
$ws.Range("N7").Value = 43.796465
$ws.Range("O7").Value = 0.2699193193490487
$ws.Range("P7").Value = 0.2699193193490487
$ws.Range("Q7").Value = 19163.92357750063
$ws.Range("R7").Value = 172475.3121975057
$ws.Range("S7").Value = 0.2607045169732969
$ws.Range("T7").Value = 0.2607045169732969
$ws.Range("G8").Value = 45.356022
$ws.Range("H8").Value = 136.068066
$ws.Range("I8").Value = 0.03337205234791334
$ws.Range("J8").Value = 0.03337205234791334
$ws.Range("M8").Value = 7.487621999999999
$ws.Range("N8").Value = 22.462866
$ws.Range("O8").Value = 0.1384395179233961
$ws.Range("P8").Value = 0.1384395179233961
$ws.Range("Q8").Value = 339.608748159684
$ws.Range("R8").Value = 3056.478733437155
$ws.Range("S8").Value = 0.00462001083915946
$ws.Range("T8").Value = 0.00462001083915946
$ws.Range("G9").Value = 45.356022
$ws.Range("H9").Value = 136.068066
$ws.Range("I9").Value = 0.03337205234791334
$ws.Range("J9").Value = 0.03337205234791334
$ws.Range("O9").Value = 0.5916411627275552
$ws.Range("P9").Value = 0.5916411627275552
$ws.Range("Q9").Value = 1451.366760355416
$ws.Range("R9").Value = 13062.30084319874
$ws.Range("S9").Value = 0.01974427985372429
$ws.Range("T9").Value = 0.01974427985372429
$ws.Range("G10").Value = 45.356022
$ws.Range("H10").Value = 136.068066
$ws.Range("I10").Value = 0.03337205234791334
$ws.Range("J10").Value = 0.03337205234791334
$ws.Range("M10").Value = 14.59882166666667
$ws.Range("N10").Value = 43.796465
$ws.Range("O10").Value = 0.2699193193490487
$ws.Range("P10").Value = 0.2699193193490487
$ws.Range("Q10").Value = 662.1444766874099
$ws.Range("R10").Value = 5959.300290186689
$ws.Range("S10").Value = 0.009007761655029591
$ws.Range("T10").Value = 0.009007761655029591
